$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: collapse the nine separate runs ("Applications", " ",
#    "of", " ", "Gaussian", " ", "Elimination", " ", "Questions") into a
#    single run holding the full text, exactly as the target XML does.
# ---------------------------------------------------------------------------
$titlePar = $d.Paragraphs(1)
$titleRng = $titlePar.Range
$titleRng.MoveEnd(1, -1) | Out-Null          # exclude the paragraph mark
$titleText = $titleRng.Text                   # runs concatenate -> full title
$titleXml = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r>' +
            '<w:t xml:space="preserve">' + $titleText + '</w:t></w:r></w:p>'
$titleRng.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) Every matrix equation's enclosing delimiter ("[ ... ]") declares its
#    characters via <m:dPr><m:begChr/><m:endChr/><m:sepChr/><m:grow/></m:dPr>.
#    Swap the serialized order of <m:endChr/> and <m:sepChr/> so it reads
#    begChr, sepChr, endChr, grow - without touching any other markup or any
#    of the matrix content itself.
# ---------------------------------------------------------------------------
$full = $d.Content.WordOpenXML

$partMarker = '<pkg:part pkg:name="/word/document.xml"'
$partIdx = $full.IndexOf($partMarker)
$dataStart = $full.IndexOf('<pkg:xmlData>', $partIdx) + '<pkg:xmlData>'.Length
$dataEnd = $full.IndexOf('</pkg:xmlData>', $dataStart)
$docXml = $full.Substring($dataStart, $dataEnd - $dataStart)

$paraPattern = [regex]'<m:oMathPara>.*?</m:oMathPara>'
$paraMatches = $paraPattern.Matches($docXml)

$oldOrder = '<m:begChr m:val="[" /><m:endChr m:val="]" /><m:sepChr m:val="" />'
$newOrder = '<m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" />'

$fixedBlocks = @()
foreach ($m in $paraMatches) {
    $block = $m.Value
    if ($block.Contains($oldOrder)) {
        $block = $block.Replace($oldOrder, $newOrder)
    }
    $fixedBlocks += , $block
}

for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $rng = $om.Range
    $rng.InsertXML($fixedBlocks[$i - 1])
}

Write-Output "edit complete"
